# ST2 Tool Monitoring sheet rework:
#  - rename header columns to the new tool-monitoring checklist naming
#  - insert 4 new tracked columns (U-X) for the "Spring Seat Assembly_BF_..."
#    checks, which pushes Remark / QA-Sign / Engg-Sign from U:W out to Y:AA
#  - update the sample data row (row 2) accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: relocate the trailing fixed columns (Remark, QA-Sign, Engg-Sign)
# from U:W to Y:AA *before* anything else is overwritten - single Copy(Destination)
# call moves both values and formatting (keeps header style s="1" on row 1,
# default style on row 2) in one shot.
$ws.Range("U1:W2").Copy($ws.Range("Y1:AA2"))
$excel.CutCopyMode = $false

# --- Step 2: header row (row 1) text updates -------------------------------
$ws.Range("C1").Value = "Lip Ring Assembly_TF_Check Fixture_OBS"
$ws.Range("D1").Value = "Lip Ring Assembly_TF_Check Fixture_Tool_Life_cnt"
$ws.Range("E1").Value = "Lip Ring Assembly_TF_Check the Tip of the Fixture_OBS"
$ws.Range("F1").Value = "Lip Ring Assembly_TF_Check the Tip of the Fixture_Tool_Life_cnt"
$ws.Range("G1").Value = "Lip Ring Assembly_TF_Check Fixture wear out_OBS"
$ws.Range("H1").Value = "Lip Ring Assembly_TF_Check Fixture wear out_Tool_Life_cnt"
$ws.Range("I1").Value = "Lip Ring Assembly_BF_Check Fixture_OBS"
$ws.Range("J1").Value = "Lip Ring Assembly_BF_Check Fixture_Tool_Life_cnt"
$ws.Range("K1").Value = "Lip Ring Assembly_BF_Check Fixture wear out_OBS"
$ws.Range("L1").Value = "Lip Ring Assembly_BF_Check Fixture wear out_Tool_Life_cnt"
$ws.Range("M1").Value = "Spring Seat Assembly_TF_Check Fixture_OBS"
$ws.Range("N1").Value = "Spring Seat Assembly_TF_Check Fixture_Tool_Life_cnt"
$ws.Range("O1").Value = "Spring Seat Assembly_TF_Check the Tip of the Fixture_OBS"
$ws.Range("P1").Value = "Spring Seat Assembly_TF_Check the Tip of the Fixture_Tool_Life_cnt"
$ws.Range("Q1").Value = "Spring Seat Assembly_TF_Check Fixture wear out_OBS"
$ws.Range("R1").Value = "Spring Seat Assembly_TF_Check Fixture wear out_Tool_Life_cnt"
$ws.Range("S1").Value = "Spring Seat Assembly_BF_Check Fixture_OBS"
$ws.Range("T1").Value = "Spring Seat Assembly_BF_Check Fixture_Tool_Life_cnt"
$ws.Range("U1").Value = "Spring Seat Assembly_BF_Check the Tip of the Fixture_OBS"
$ws.Range("V1").Value = "Spring Seat Assembly_BF_Check the Tip of the Fixture_Tool_Life_cnt"
$ws.Range("W1").Value = "Spring Seat Assembly_BF_Check Fixture wear out_OBS"
$ws.Range("X1").Value = "Spring Seat Assembly_BF_Check Fixture wear out_Tool_Life_cnt"

# New header cells (U1:X1) need the same bold/border/centered style as the
# rest of the header row - copy format+value from an existing header cell
# onto them (values were already set above, so re-apply just the texts after).
$ws.Range("T1").Copy($ws.Range("U1:X1"))
$excel.CutCopyMode = $false
$ws.Range("U1").Value = "Spring Seat Assembly_BF_Check the Tip of the Fixture_OBS"
$ws.Range("V1").Value = "Spring Seat Assembly_BF_Check the Tip of the Fixture_Tool_Life_cnt"
$ws.Range("W1").Value = "Spring Seat Assembly_BF_Check Fixture wear out_OBS"
$ws.Range("X1").Value = "Spring Seat Assembly_BF_Check Fixture wear out_Tool_Life_cnt"

# --- Step 3: sample data row (row 2) updates -------------------------------
$ws.Range("A2").Value = "2025-02-06T14:07"
$ws.Range("B2").Value = "SHIFT2"

# numeric-looking counters are stored as text in this template, so force the
# text number-format before assigning the value to keep them as strings
$countCells = @("F2", "H2", "J2", "L2", "N2", "P2", "R2", "T2", "V2", "X2")
foreach ($addr in $countCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("F2").Value = "2"
$ws.Range("H2").Value = "3"
$ws.Range("J2").Value = "4"
$ws.Range("L2").Value = "6"
$ws.Range("N2").Value = "7"
$ws.Range("P2").Value = "8"
$ws.Range("R2").Value = "9"
$ws.Range("T2").Value = "10"
$ws.Range("V2").Value = "11"
$ws.Range("X2").Value = "12"

$ws.Range("U2").Value = "OK"
$ws.Range("W2").Value = "OK"

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("Y2").Value = "000"
$ws.Range("Z2").Value = "111"
$ws.Range("AA2").Value = "222"

$wb.Save()
